# Weekly refresh of fruit/vegetable price data ("Mora" subset).
# The D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen) and S (Precio $/Kg) columns
# are updated per row 2-17 to reflect the new weekly snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
  2  = @(44533, 150, 4000, 4000, 4000, "Provincia de Curicó", 2000)
  3  = @(44231, 150, 3400, 3400, 3400, "Provincia de Curicó", 1700)
  4  = @(44237, 100, 3600, 4000, 3800, "Provincia de Curicó", 1900)
  5  = @(44188, 150, 3000, 3400, 3240, "Provincia de Linares", 1620)
  6  = @(44236, 300, 3600, 4000, 3800, "Provincia de Curicó", 1900)
  7  = @(44168, 170, 8000, 8000, 8000, "Provincia de Linares", 4000)
  8  = @(44586, 250, 5000, 5000, 5000, "Provincia de Curicó", 2500)
  9  = @(44238, 300, 3600, 4000, 3800, "Provincia de Curicó", 1900)
  10 = @(44582, 380, 5000, 5000, 5000, "Provincia de Curicó", 2500)
  11 = @(44980, 250, 4000, 4000, 4000, "Provincia de Curicó", 2000)
  12 = @(44194, 120, 3000, 3000, 3000, "Provincia de Linares", 1500)
  13 = @(44617, 90,  6500, 6500, 6500, "Provincia de Curicó", 3250)
  14 = @(44208, 85,  3000, 3000, 3000, "Provincia de Linares", 1500)
  15 = @(44174, 200, 3200, 3200, 3200, "Provincia de Curicó", 1600)
  16 = @(44978, 500, 3000, 3000, 3000, "Provincia de Curicó", 1500)
  17 = @(44232, 200, 3000, 3000, 3000, "Provincia de Curicó", 1500)
}

foreach ($r in $rows.Keys) {
  $v = $rows[$r]
  $ws.Cells.Item($r, 4).Value  = $v[0]   # D: Fecha
  $ws.Cells.Item($r, 13).Value = $v[1]   # M: Volumen
  $ws.Cells.Item($r, 14).Value = $v[2]   # N: Precio minimo
  $ws.Cells.Item($r, 15).Value = $v[3]   # O: Precio maximo
  $ws.Cells.Item($r, 16).Value = $v[4]   # P: Precio promedio ponderado
  $ws.Cells.Item($r, 18).Value = $v[5]   # R: Origen
  $ws.Cells.Item($r, 19).Value = $v[6]   # S: Precio $/Kg
}
